$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they stay plain text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.915.35'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '3.037.63'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '594.88'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').Value = '153.69'
$ws.Range('E6').Value = '  +7.43%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.034.61'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '6.84'
$ws.Range('E10').Value = '  +13.87%  '
$ws.Range('E11').Value = '  +4.42%  '
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').Value = '  +2.66%  '
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('D14').Value = '35.81'
$ws.Range('E14').Value = '  +4.74%  '
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '3.544.95'
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').Value = '7.09'
$ws.Range('E17').Value = '  +3.00%  '
$ws.Range('D18').Value = '62.922.93'
$ws.Range('E18').Value = '  +2.70%  '
$ws.Range('D19').Value = '3.041.53'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').Value = '453.78'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').Value = '14.30'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('D22').Value = '0.698'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').Value = '7.53'
$ws.Range('E23').Value = '  +3.37%  '
$ws.Range('D24').Value = '83.07'
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('D25').Value = '11.25'
$ws.Range('E25').Value = '  +6.44%  '
$ws.Range('D26').Value = '2.30'
$ws.Range('E26').Value = '  +5.08%  '
$ws.Range('D27').Value = '12.45'
$ws.Range('E27').Value = '  +4.81%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '7.45'
$ws.Range('E29').Value = '  +4.91%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '2.25'
$ws.Range('E30').Value = '  +10.41%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').Value = '27.63'
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').Value = '0.0₃0857'
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D36').Value = '1.05'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('D37').Value = '5.94'
$ws.Range('E37').Value = '  +3.42%  '
$ws.Range('D38').Value = '3.17'
$ws.Range('E38').Value = '  +10.32%  '
$ws.Range('D39').Value = '0.131'
$ws.Range('E39').Value = '  +8.28%  '
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('D41').Value = '50.39'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').Value = '0.305'
$ws.Range('E43').Value = '  +13.97%  '
$ws.Range('D44').Value = '43.75'
$ws.Range('E44').Value = '  +10.45%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0363'
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '391.35'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '2.722.08'
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('D48').Value = '132.65'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +7.67%  '
$ws.Range('D51').Value = '24.55'
$ws.Range('E51').Value = '  +5.23%  '
